$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99 (shifts existing rows 99-154 down to 100-155,
# carrying their original values/formatting with them).
$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new weekly price record.
$ws.Range("A99").Value = 1
$ws.Range("B99").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C99").Value = "Arica y Parinacota"
$ws.Range("D99").Value = 44942
$ws.Range("E99").Value = 15
$ws.Range("F99").Value = 100112042
$ws.Range("G99").Value = "Locoto"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 90
$ws.Range("K99").Value = 19000
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = 19667
$ws.Range("N99").Value = "$/caja 20 kilos"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 983
$ws.Range("Q99").Value = 20
$ws.Range("R99").Value = "Hortaliza"

# D99 needs the same date-number-format style as the rest of the "Fecha" column.
$ws.Range("D99").NumberFormat = $ws.Range("D100").NumberFormat
